$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 8 ("Alors, on commence !"),
# pushing it (and everything below) down to rows 10+.
$ws.Rows("8:9").Insert()

# New row 8 - first new instruction slide
$ws.Range("A8").Value = "Instruct"
$ws.Range("B8").Value = "Dans ce qui suit tu vas bien entendre les sauts`net essayer d'appuyer sur la bonne touche."
$ws.Range("C8").Value = "Diapositive6b"
$ws.Range("D8").Value = "Key"
$ws.Range("E8").Value = "None"

# New row 9 - second new instruction slide
$ws.Range("A9").Value = "Instruct"
$ws.Range("B9").Value = "Si la première fois tu n'as pas bien entendu les sauts, tu peux réécouter`nen appuyant sur ESPACE"
$ws.Range("C9").Value = "Diapositive6c"
$ws.Range("D9").Value = "Key"
$ws.Range("E9").Value = "None"

# Apply the new paragraph-style formatting (12pt black Calibri, centered
# vertically, wrapped) used for these two new instruction cells.
$newTextRange = $ws.Range("B8:B9")
$newTextRange.Font.Size = 12
$newTextRange.Font.Color = 0
$newTextRange.WrapText = $true
$newTextRange.VerticalAlignment = -4108

# The row that used to be row 8 ("Alors, on commence !") is now row 10;
# give its text cell the new number-format + wrap style.
$ws.Range("B10").WrapText = $true
$ws.Range("B10").NumberFormat = "h:mm"

# Selection moves onto the (shifted) "Alors, on commence !" cell.
$ws.Range("B10").Select()

# Page now prints in portrait orientation.
$ws.PageSetup.Orientation = 1
